$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for columns E:T, rows 2-7
$data = @{
    2 = @("1","0.3333333333333333","0.03644533333333334","0.109336","0.005561955322140003","0.005561955322140003","3","1","0.1949313333333333","0.5847939999999999","0.07226189720071244","0.07226189720071244","0.007104337420444445","0.06393903678399999","0.0004019174437234363","0.0004019174437234363")
    3 = @("1","0.3333333333333333","0.03644533333333334","0.109336","0.005561955322140003","0.005561955322140003","3","1","2.102325333333333","6.306976000000001","0.779341189135594","0.7793411891355941","0.07661994754844445","0.6895795279360001","0.004334660874675635","0.004334660874675636")
    4 = @("1","0.3333333333333333","0.03644533333333334","0.109336","0.005561955322140003","0.005561955322140003","2","0.6666666666666666","0.4003106666666666","1.200932","0.1483969136636935","0.1483969136636935","0.01458945568355556","0.131305101152","0.0008253770037409307","0.0008253770037409307")
    5 = @("3","1","6.516166333333334","19.548499","0.99443804467786","0.99443804467786","3","1","0.1949313333333333","0.5847939999999999","0.07226189720071244","0.07226189720071244","1.270204991578444","11.431844924206","0.071859979756989","0.071859979756989")
    6 = @("3","1","6.516166333333334","19.548499","0.99443804467786","0.99443804467786","3","1","2.102325333333333","6.306976000000001","0.779341189135594","0.7793411891355941","13.69910155878044","123.291914029024","0.7750065282609184","0.7750065282609185")
    7 = @("3","1","6.516166333333334","19.548499","0.99443804467786","0.99443804467786","2","0.6666666666666666","0.4003106666666666","1.200932","0.1483969136636935","0.1483969136636935","2.608490889007555","23.476418001068","0.1475715366599526","0.1475715366599526")
}

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = "$($cols[$i])$row"
        $ws.Range($cell).Value = [double]$vals[$i]
    }
}
